$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.480.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.097.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.097.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.631.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.564.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.100.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "337.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.12%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0911"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0659"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.138.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.682"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.10%  "
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.296.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.978"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.58%  "
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E51").Value = "  +1.28%  "
